# Generate Report for Handoff
# Updates the localization-status report: the file has moved from
# "Handed back: in sync with en-US" to "Ready for handoff", with refreshed
# timestamps, across the Overview sheet and the per-language detail sheets.
# Because the status text got shorter, the Status columns are re-sized to
# fit the new content.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

$newStatus = "Ready for handoff"

# --- Overview sheet -------------------------------------------------
$wsOverview.Range("E2").Value = $newStatus
$wsOverview.Range("F2").Value = $newStatus
$wsOverview.Range("G2").Value = "2016-08-14 17:19:28"

# --- zh-cn detail sheet ----------------------------------------------
$wsZhCn.Range("C2").Value = $newStatus
$wsZhCn.Range("H2").Value = "2016-08-14 17:19:20"

# --- de-de detail sheet ----------------------------------------------
$wsDeDe.Range("C2").Value = $newStatus
$wsDeDe.Range("H2").Value = "2016-08-14 17:19:28"

# --- Resize the Status columns to fit the shorter text ----------------
# (matches the narrower column width produced when the report was
# regenerated with the new, shorter status text)
$wsOverview.Columns.Item(5).ColumnWidth = 16.33
$wsOverview.Columns.Item(6).ColumnWidth = 16.33
$wsZhCn.Columns.Item(3).ColumnWidth = 16.33
$wsDeDe.Columns.Item(3).ColumnWidth = 16.33
